$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column J (2022 data) gets corrected values ---
# Rows 2-12: new values, and formatting reverts to the default/general style
# (no explicit number format anymore).
$newValues = @{
    2  = 299226.13
    3  = 274821.67
    4  = 305662.04
    5  = 275167.42
    6  = 276489.61
    7  = 270059.27
    8  = 265810.43
    9  = 280041.17
    10 = 259983.59
    11 = 283615.59
    12 = 277227.8
}

foreach ($row in $newValues.Keys) {
    $cell = $ws.Cells.Item($row, 10)
    $cell.Style = "Normal"
    $cell.Value = $newValues[$row]
}

# Row 13 also gets a corrected value, but keeps an explicit (font-applied) style.
$j13 = $ws.Cells.Item(13, 10)
$j13.Style = "Normal"
$j13.Font.Name = "Calibri"
$j13.Value = 306083.39

# --- Page setup: A4, portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Final selection left on J13 ---
$ws.Range("J13").Select()
